$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---
$ws1.Range("H2").Value = 1.61
$ws1.Range("L2").Value = 1.11

$ws1.Range("H3").Value = 0.61
$ws1.Range("J3").Value = "Urgent"
$ws1.Range("L3").Value = 1.09

$ws1.Range("L4").Value = 1.02
$ws1.Range("L5").Value = 1.17
$ws1.Range("L6").Value = 1.12
$ws1.Range("L7").Value = 1.15
$ws1.Range("L8").Value = 1.05
$ws1.Range("L9").Value = 0.9399999999999999
$ws1.Range("L10").Value = 1.13
$ws1.Range("L11").Value = 1.14
$ws1.Range("L12").Value = 1.1
$ws1.Range("L13").Value = 0.9
$ws1.Range("L14").Value = 1.07
$ws1.Range("L15").Value = 0.85
$ws1.Range("L16").Value = 1.18
$ws1.Range("L17").Value = 1.02

# --- Summary sheet ---
# Force these as text (not numbers) since they are numeric-looking
# string values in the original inline-string cells.
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "5"
$ws2.Range("B14").NumberFormat = "@"
$ws2.Range("B14").Value = "1"
